$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 144.783305
$ws.Range("H2").Value = 434.349915
$ws.Range("I2").Value = 0.2430046335191003
$ws.Range("J2").Value = 0.251012682214973
$ws.Range("M2").Value = 1.495937333333333
$ws.Range("N2").Value = 4.487812
$ws.Range("O2").Value = 0.3639577964822184
$ws.Range("P2").Value = 0.3715212096336638
$ws.Range("Q2").Value = 216.5867511928866
$ws.Range("R2").Value = 1949.28076073598
$ws.Range("S2").Value = 0.08844343095058078
$ws.Range("T2").Value = 0.09325653532989719

$ws.Range("G3").Value = 144.783305
$ws.Range("H3").Value = 434.349915
$ws.Range("I3").Value = 0.2430046335191003
$ws.Range("J3").Value = 0.251012682214973
$ws.Range("O3").Value = 0.5297630075657004
$ws.Range("P3").Value = 0.5407720216252946
$ws.Range("Q3").Value = 315.255367022845
$ws.Range("R3").Value = 2837.298303205605
$ws.Range("S3").Value = 0.1287348655054794
$ws.Range("T3").Value = 0.1357406356149786

$ws.Range("G4").Value = 144.783305
$ws.Range("H4").Value = 434.349915
$ws.Range("I4").Value = 0.2430046335191003
$ws.Range("J4").Value = 0.251012682214973
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.2510255
$ws.Range("N4").Value = 0.502051
$ws.Range("O4").Value = 0.06107387375463619
$ws.Range("P4").Value = 0.04156203397508419
$ws.Range("Q4").Value = 36.34430152927749
$ws.Range("R4").Value = 218.065809175665
$ws.Range("S4").Value = 0.01484123430933717
$ws.Range("T4").Value = 0.01043259762639572

$ws.Range("G5").Value = 144.783305
$ws.Range("H5").Value = 434.349915
$ws.Range("I5").Value = 0.2430046335191003
$ws.Range("J5").Value = 0.251012682214973
$ws.Range("M5").Value = 0.1858026666666667
$ws.Range("N5").Value = 0.557408
$ws.Range("O5").Value = 0.04520532219744508
$ws.Range("P5").Value = 0.0461447347659575
$ws.Range("Q5").Value = 26.90112415781333
$ws.Range("R5").Value = 242.11011742032
$ws.Range("S5").Value = 0.01098510275370299
$ws.Range("T5").Value = 0.01158291364370151

$ws.Range("G6").Value = 82.24887099999999
$ws.Range("I6").Value = 0.1380466950572427
$ws.Range("J6").Value = 0.1425959278859072
$ws.Range("M6").Value = 1.495937333333333
$ws.Range("N6").Value = 4.487812
$ws.Range("O6").Value = 0.3639577964822184
$ws.Range("P6").Value = 0.3715212096336638
$ws.Range("Q6").Value = 123.0391567534173
$ws.Range("R6").Value = 1107.352410780756
$ws.Range("S6").Value = 0.05024317094468679
$ws.Range("T6").Value = 0.05297741161700693

$ws.Range("G7").Value = 82.24887099999999
$ws.Range("I7").Value = 0.1380466950572427
$ws.Range("J7").Value = 0.1425959278859072
$ws.Range("O7").Value = 0.5297630075657004
$ws.Range("P7").Value = 0.5407720216252946
$ws.Range("S7").Value = 0.07313203235802999
$ws.Range("T7").Value = 0.07711188819839675

$ws.Range("G8").Value = 82.24887099999999
$ws.Range("I8").Value = 0.1380466950572427
$ws.Range("J8").Value = 0.1425959278859072
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 0.2510255
$ws.Range("N8").Value = 0.502051
$ws.Range("O8").Value = 0.06107387375463619
$ws.Range("P8").Value = 0.04156203397508419
$ws.Range("Q8").Value = 20.6465639672105
$ws.Range("R8").Value = 123.879383803263
$ws.Range("S8").Value = 0.008431046426170799
$ws.Range("T8").Value = 0.005926576799502729

$ws.Range("G9").Value = 82.24887099999999
$ws.Range("I9").Value = 0.1380466950572427
$ws.Range("J9").Value = 0.1425959278859072
$ws.Range("M9").Value = 0.1858026666666667
$ws.Range("N9").Value = 0.557408
$ws.Range("O9").Value = 0.04520532219744508
$ws.Range("P9").Value = 0.0461447347659575
$ws.Range("Q9").Value = 15.28205956212267
$ws.Range("R9").Value = 137.538536059104
$ws.Range("S9").Value = 0.006240445328355104
$ws.Range("T9").Value = 0.00658005127100079

$ws.Range("G10").Value = 163.8590903333333
$ws.Range("H10").Value = 491.577271
$ws.Range("I10").Value = 0.2750214756820535
$ws.Range("J10").Value = 0.284084617144743
$ws.Range("M10").Value = 1.495937333333333
$ws.Range("N10").Value = 4.487812
$ws.Range("O10").Value = 0.3639577964822184
$ws.Range("P10").Value = 0.3715212096336638
$ws.Range("Q10").Value = 245.1229306356724
$ws.Range("R10").Value = 2206.106375721052
$ws.Range("S10").Value = 0.1000962102745282
$ws.Range("T10").Value = 0.1055434605999312

$ws.Range("G11").Value = 163.8590903333333
$ws.Range("H11").Value = 491.577271
$ws.Range("I11").Value = 0.2750214756820535
$ws.Range("J11").Value = 0.284084617144743
$ws.Range("O11").Value = 0.5297630075657004
$ws.Range("P11").Value = 0.5407720216252946
$ws.Range("Q11").Value = 356.7915352054197
$ws.Range("R11").Value = 3211.123816848777
$ws.Range("S11").Value = 0.1456962041024818
$ws.Range("T11").Value = 0.1536250127260105

$ws.Range("G12").Value = 163.8590903333333
$ws.Range("H12").Value = 491.577271
$ws.Range("I12").Value = 0.2750214756820535
$ws.Range("J12").Value = 0.284084617144743
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 0.2510255
$ws.Range("N12").Value = 0.502051
$ws.Range("O12").Value = 0.06107387375463619
$ws.Range("P12").Value = 0.04156203397508419
$ws.Range("Q12").Value = 41.13281008047017
$ws.Range("R12").Value = 246.796860482821
$ws.Range("S12").Value = 0.01679662688561948
$ws.Range("T12").Value = 0.01180713450956859

$ws.Range("G13").Value = 163.8590903333333
$ws.Range("H13").Value = 491.577271
$ws.Range("I13").Value = 0.2750214756820535
$ws.Range("J13").Value = 0.284084617144743
$ws.Range("M13").Value = 0.1858026666666667
$ws.Range("N13").Value = 0.557408
$ws.Range("O13").Value = 0.04520532219744508
$ws.Range("P13").Value = 0.0461447347659575
$ws.Range("Q13").Value = 30.44545594150756
$ws.Range("R13").Value = 274.009103473568
$ws.Range("S13").Value = 0.01243243441942403
$ws.Range("T13").Value = 0.01310900930923275

$ws.Range("G14").Value = 57.0238095
$ws.Range("H14").Value = 114.047619
$ws.Range("I14").Value = 0.09570889357312636
$ws.Range("J14").Value = 0.06590860906562239
$ws.Range("M14").Value = 1.495937333333333
$ws.Range("N14").Value = 4.487812
$ws.Range("O14").Value = 0.3639577964822184
$ws.Range("P14").Value = 0.3715212096336638
$ws.Range("Q14").Value = 85.30404551993799
$ws.Range("R14").Value = 511.824273119628
$ws.Range("S14").Value = 0.03483399800862622
$ws.Range("T14").Value = 0.02448644616533229

$ws.Range("G15").Value = 57.0238095
$ws.Range("H15").Value = 114.047619
$ws.Range("I15").Value = 0.09570889357312636
$ws.Range("J15").Value = 0.06590860906562239
$ws.Range("O15").Value = 0.5297630075657004
$ws.Range("P15").Value = 0.5407720216252946
$ws.Range("Q15").Value = 124.1652964957755
$ws.Range("R15").Value = 744.991778974653
$ws.Range("S15").Value = 0.05070303131008495
$ws.Range("T15").Value = 0.03564153176692784

$ws.Range("G16").Value = 57.0238095
$ws.Range("H16").Value = 114.047619
$ws.Range("I16").Value = 0.09570889357312636
$ws.Range("J16").Value = 0.06590860906562239
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.5
$ws.Range("M16").Value = 0.2510255
$ws.Range("N16").Value = 0.502051
$ws.Range("O16").Value = 0.06107387375463619
$ws.Range("P16").Value = 0.04156203397508419
$ws.Range("Q16").Value = 14.31443029164225
$ws.Range("R16").Value = 57.257721166569
$ws.Range("S16").Value = 0.005845312883281031
$ws.Range("T16").Value = 0.002739295849235939

$ws.Range("G17").Value = 57.0238095
$ws.Range("H17").Value = 114.047619
$ws.Range("I17").Value = 0.09570889357312636
$ws.Range("J17").Value = 0.06590860906562239
$ws.Range("M17").Value = 0.1858026666666667
$ws.Range("N17").Value = 0.557408
$ws.Range("O17").Value = 0.04520532219744508
$ws.Range("P17").Value = 0.0461447347659575
$ws.Range("Q17").Value = 10.595175868592
$ws.Range("R17").Value = 63.571055211552
$ws.Range("S17").Value = 0.004326551371134158
$ws.Range("T17").Value = 0.003041335284126327

$ws.Range("G18").Value = 147.8896333333333
$ws.Range("H18").Value = 443.6689
$ws.Range("I18").Value = 0.2482183021684772
$ws.Range("J18").Value = 0.2563981636887546
$ws.Range("M18").Value = 1.495937333333333
$ws.Range("N18").Value = 4.487812
$ws.Range("O18").Value = 0.3639577964822184
$ws.Range("P18").Value = 0.3715212096336638
$ws.Range("Q18").Value = 221.2336237163111
$ws.Range("R18").Value = 1991.1026134468
$ws.Range("S18").Value = 0.09034098630379642
$ws.Range("T18").Value = 0.09525735592149621

$ws.Range("G19").Value = 147.8896333333333
$ws.Range("H19").Value = 443.6689
$ws.Range("I19").Value = 0.2482183021684772
$ws.Range("J19").Value = 0.2563981636887546
$ws.Range("O19").Value = 0.5297630075657004
$ws.Range("P19").Value = 0.5407720216252946
$ws.Range("Q19").Value = 322.0191764193667
$ws.Range("R19").Value = 2898.1725877743
$ws.Range("S19").Value = 0.1314968742896243
$ws.Range("T19").Value = 0.138652953318981

$ws.Range("G20").Value = 147.8896333333333
$ws.Range("H20").Value = 443.6689
$ws.Range("I20").Value = 0.2482183021684772
$ws.Range("J20").Value = 0.2563981636887546
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.5
$ws.Range("M20").Value = 0.2510255
$ws.Range("N20").Value = 0.502051
$ws.Range("O20").Value = 0.06107387375463619
$ws.Range("P20").Value = 0.04156203397508419
$ws.Range("Q20").Value = 37.12406915231667
$ws.Range("R20").Value = 222.7444149139
$ws.Range("S20").Value = 0.01515965325022772
$ws.Range("T20").Value = 0.01065642919038121

$ws.Range("G21").Value = 147.8896333333333
$ws.Range("H21").Value = 443.6689
$ws.Range("I21").Value = 0.2482183021684772
$ws.Range("J21").Value = 0.2563981636887546
$ws.Range("M21").Value = 0.1858026666666667
$ws.Range("N21").Value = 0.557408
$ws.Range("O21").Value = 0.04520532219744508
$ws.Range("P21").Value = 0.0461447347659575
$ws.Range("Q21").Value = 27.47828824568889
$ws.Range("R21").Value = 247.3045942112
$ws.Range("S21").Value = 0.01122078832482879
$ws.Range("T21").Value = 0.01183142525789613
